# Asset list now imports RFID tags: replace the old "AssetID/Type" table
# with the new "Asset ID/RFID Tag" table (4 data rows instead of 11) and
# tidy up the sheet layout (column B width, selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Asset ID"
$ws.Range("B1").Value = "RFID Tag"

# --- Data rows (Asset ID -> RFID Tag) --------------------------------
$ws.Range("A2").Value = "E1010101"
$ws.Range("B2").Value = "45345345zxvzxvzxb"

$ws.Range("A3").Value = "E0101010"
$ws.Range("B3").Value = "45245345345asrasrsar"

$ws.Range("A4").Value = "E6996696"
$ws.Range("B4").Value = "waweaser"

$ws.Range("A5").Value = "4867530"
$ws.Range("B5").Value = "4524534534werwerwerwer"

# Remove the now-unused rows 6-12 that held the old asset list, shrinking
# the sheet's used range down to A1:B5.
$ws.Range("A6:A12").EntireRow.Delete()

# New RFID Tag column needs to be wide enough to show its values.
$ws.Columns("B:B").ColumnWidth = 11.14

# Leave the selection on the last data row, column A, like the author did.
[void]$ws.Range("A5").Select()
